$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("ورودی واگن یا بار"): update a few numeric inputs ---
$ws1.Range("B4").Value = 1350
$ws1.Range("A7").Value = 2400
$ws1.Range("A8").Value = 3000

# --- Sheet2 ("خروجی"): add a "فضای سازه" (structure-space) column and
#     rework the passability wording / values for the gabari summary ---

# Insert a brand-new column G (inherits each row's formatting from the
# neighboring column F, so every new cell reuses the existing per-row
# style index instead of creating new styles.xml entries).
$ws2.Columns("G:G").Insert()
$ws2.Columns("G:G").ColumnWidth = 20.166666666666668

# Only G1 and G6 actually hold data in the final layout - drop the
# placeholder cells Insert() created for rows 2-5 so they disappear
# entirely (matching cells A2:F... which never had a G neighbor).
$ws2.Range("G2:G5").Clear()

# Row 1 (headers)
$ws2.Range("C1").Value = "قابلیت عبور از فضای مجاز"
$ws2.Range("E1").Value = "قابلیت عبور از فضای آزاد"
$ws2.Range("F1").Value = "اندازه ورود به فضای سازه"
$ws2.Range("G1").Value = "قابلیت عبور از فضای سازه"

# Row 2
$ws2.Range("D2").Value = 70.71067811865476
$ws2.Range("F2").Clear()

# Row 3
$ws2.Range("C3").Value = "قابل عبور"
$ws2.Range("D3:F3").Clear()

# Row 4
$ws2.Range("D4").Value = 70.71067811865476
$ws2.Range("F4").Clear()

# Row 5
$ws2.Range("C5").Value = "قابل عبور"
$ws2.Range("D5:F5").Clear()

# Row 6
$ws2.Range("D6").Value = 340.94574348420895
$ws2.Range("F6").Value = 150
$ws2.Range("G6").Value = "غیر قابل عبور"

Write-Host "edits applied"
